$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.360.68'
$ws.Range('E2').Value = '  -0.19%  '

# Row 3
$ws.Range('D3').Value = '1.887.63'
$ws.Range('E3').Value = '  -1.54%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.78%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '340.12'
$ws.Range('E5').Value = '  +4.80%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  -0.68%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4729'
$ws.Range('E7').Value = '  -1.42%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3978'
$ws.Range('E8').Value = '  -1.71%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.11'
$ws.Range('E9').Value = '  -1.65%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07995'
$ws.Range('E10').Value = '  -2.23%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9889'
$ws.Range('E11').Value = '  -1.88%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.84'
$ws.Range('E12').Value = '  -2.33%  '

# Row 13
$ws.Range('D13').Value = '1.933.77'
$ws.Range('E13').Value = '  +0.63%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.910'
$ws.Range('E14').Value = '  -2.21%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.082'
$ws.Range('E15').Value = '  -1.95%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.88'
$ws.Range('E16').Value = '  -2.90%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.06771'
$ws.Range('E17').Value = '  -1.38%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.004'
$ws.Range('E18').Value = '  -0.63%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001015'
$ws.Range('E19').Value = '  -2.07%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.26'
$ws.Range('E20').Value = '  -1.68%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  -0.59%  '

# Row 22
$ws.Range('D22').Value = '29.376.75'
$ws.Range('E22').Value = '  -0.16%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.492'
$ws.Range('E23').Value = '  -2.95%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.53'
$ws.Range('E24').Value = '  -2.36%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.146'
$ws.Range('E25').Value = '  -2.10%  '

# Row 26
$ws.Range('D26').Value = '2.174.47'
$ws.Range('E26').Value = '  +0.98%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '155.85'
$ws.Range('E27').Value = '  +0.21%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.450'
$ws.Range('E28').Value = '  -1.86%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.60'
$ws.Range('E29').Value = '  -1.92%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.036'
$ws.Range('E30').Value = '  -3.03%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '118.58'
$ws.Range('E31').Value = '  -1.43%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9948'
$ws.Range('E32').Value = '  -1.62%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09499'
$ws.Range('E33').Value = '  -1.53%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.467'
$ws.Range('E34').Value = '  -1.94%  '

# Row 35
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.384'
$ws.Range('E35').Value = '  +0.95%  '

# Row 36
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.522'
$ws.Range('E36').Value = '  -0.89%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06254'
$ws.Range('E37').Value = '  -1.92%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02235'

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.190'
$ws.Range('E39').Value = '  +0.74%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5807'
$ws.Range('E40').Value = '  -2.10%  '

# Row 41
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.712'
$ws.Range('E41').Value = '  -3.54%  '

# Row 42
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '10.25'
$ws.Range('E42').Value = '  -4.49%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1818'
$ws.Range('E43').Value = '  -1.59%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.386'
$ws.Range('E44').Value = '  +0.72%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.269'
$ws.Range('E45').Value = '  -1.02%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.09'
$ws.Range('E46').Value = '  -2.29%  '

# Row 47
$ws.Range('E47').Value = '  -1.32%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.07277'
$ws.Range('E48').Value = '  -2.92%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.954'
$ws.Range('E49').Value = '  +0.63%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '116.49'
$ws.Range('E50').Value = '  -1.54%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.367'
$ws.Range('E51').Value = '  -2.55%  '
